$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers: BTec_Logo-Orange picture, "image1.jpg" -> "image2.jpg" ---
$hdrPrimary = $sec.Headers.Item(1)
if ($hdrPrimary.Exists) {
    $shp = $hdrPrimary.Range.InlineShapes.Item(1)
    if ($shp.Name -ne "image2.jpg") {
        $shp.Name = "image2.jpg"
    }
}

$hdrFirst = $sec.Headers.Item(2)
if ($hdrFirst.Exists) {
    $shp = $hdrFirst.Range.InlineShapes.Item(1)
    if ($shp.Name -ne "image2.jpg") {
        $shp.Name = "image2.jpg"
    }
}

# --- Footers: Pearson logo picture, "image2.png" -> "image1.png" ---
$ftrPrimary = $sec.Footers.Item(1)
if ($ftrPrimary.Exists) {
    $shp = $ftrPrimary.Range.InlineShapes.Item(1)
    if ($shp.Name -ne "image1.png") {
        $shp.Name = "image1.png"
    }
}

$ftrFirst = $sec.Footers.Item(2)
if ($ftrFirst.Exists) {
    $shp = $ftrFirst.Range.InlineShapes.Item(1)
    if ($shp.Name -ne "image1.png") {
        $shp.Name = "image1.png"
    }
}
